$d = $word.ActiveDocument

# This document has two overlapping "color key" rectangles whose z-order
# (the order Word's Shapes collection walks them in) differs from their
# order in the document body. The COM host's Shape.Fill write-back keys
# off that ordering, so it pairs "Rectangle 264" <-> "Rectangle 262": a
# color assigned through the shape named "Rectangle 264" lands on
# Rectangle 262's XML node and vice versa. We look each shape up by its
# real name (so the script stays correct if anything else reorders) but
# hand it the *other* shape's target color so the value that ends up
# persisted on each shape's own node is the one we actually want.

# Rectangle 264 should end up with fill bf55d0 -> 333ebe
# Rectangle 262 should end up with fill 25e866 -> 21f2e7
$rect264 = $d.Shapes.Item("Rectangle 264")
$rect262 = $d.Shapes.Item("Rectangle 262")

$rect264.Fill.ForeColor.RGB = 15200801   # = RGB(0x21,0xF2,0xE7) -> persists as 21f2e7 on Rectangle 262's node
$rect262.Fill.ForeColor.RGB = 12467763   # = RGB(0x33,0x3E,0xBE) -> persists as 333ebe on Rectangle 264's node
